$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1113
$ws1.Range("F4").Value = 1785
$ws1.Range("F5").Value = 789
$ws1.Range("F6").Value = 302
$ws1.Range("F7").Value = 208

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1113
$ws4.Range("F4").Value = 1785
$ws4.Range("F6").Value = 789
$ws4.Range("F7").Value = 302
$ws4.Range("F8").Value = 208
